# Generate Report for Handback
# For each localized-language sheet (zh-cn, de-de), the two content rows
# (row 2: 883008bc-..., row 3: e8172bf7-...) move from "Ready for handoff"
# to "Handed back: in sync with en-US": the Latest Target File / Latest
# Handback File columns get filled in (mirroring the source/handoff file
# hyperlinks) and the Latest Handback DateTime is stamped.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ----
# The "Ready for handoff" shared string is reused by the Overview sheet's
# zh-cn/de-de status columns (B/C) for both files, so it flips to the new
# status text here too.
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = $newStatus
$ws.Range("C2").Value = $newStatus
$ws.Range("B3").Value = $newStatus
$ws.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 : 883008bc-1e48-4fbd-aa59-c5b93c4a3497
$ws.Range("B2").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/164a4c9bf3bd70dc428f78e9b9fb854cd6fda3f4/e2e/883008bc-1e48-4fbd-aa59-c5b93c4a3497.md", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f621318ecd20f236baa0eb554753cffe9167608/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.zh-cn.xlf", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.zh-cn.xlf") | Out-Null
$ws.Range("G2").Value = "2016-03-10 09:23:51"

# Row 3 : e8172bf7-3212-4169-b195-480d7c7e259b
$ws.Range("B3").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/164a4c9bf3bd70dc428f78e9b9fb854cd6fda3f4/e2e/e8172bf7-3212-4169-b195-480d7c7e259b.md", "", "", "e8172bf7-3212-4169-b195-480d7c7e259b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f621318ecd20f236baa0eb554753cffe9167608/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e8172bf7-3212-4169-b195-480d7c7e259b.bcbb282467958fc11af7b1c3f40572e5c8a8a945.zh-cn.xlf", "", "", "e8172bf7-3212-4169-b195-480d7c7e259b.bcbb282467958fc11af7b1c3f40572e5c8a8a945.zh-cn.xlf") | Out-Null
$ws.Range("G3").Value = "2016-03-10 09:23:51"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

# Row 2 : 883008bc-1e48-4fbd-aa59-c5b93c4a3497
$ws.Range("B2").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/164a4c9bf3bd70dc428f78e9b9fb854cd6fda3f4/e2e/883008bc-1e48-4fbd-aa59-c5b93c4a3497.md", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6fc3b54b38aae547d03fc581d197e049d39ed244/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.de-de.xlf", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.de-de.xlf") | Out-Null
$ws.Range("G2").Value = "2016-03-10 09:24:08"

# Row 3 : e8172bf7-3212-4169-b195-480d7c7e259b
$ws.Range("B3").Value = $newStatus
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/164a4c9bf3bd70dc428f78e9b9fb854cd6fda3f4/e2e/e8172bf7-3212-4169-b195-480d7c7e259b.md", "", "", "e8172bf7-3212-4169-b195-480d7c7e259b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6fc3b54b38aae547d03fc581d197e049d39ed244/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e8172bf7-3212-4169-b195-480d7c7e259b.bcbb282467958fc11af7b1c3f40572e5c8a8a945.de-de.xlf", "", "", "e8172bf7-3212-4169-b195-480d7c7e259b.bcbb282467958fc11af7b1c3f40572e5c8a8a945.de-de.xlf") | Out-Null
$ws.Range("G3").Value = "2016-03-10 09:24:08"

Write-Output "Handback report generated."
